$wb = $excel.ActiveWorkbook

# --- ModuleName sheet: add new row "Engagements", resize column A, and make this the active sheet/tab ---
$wsModule = $wb.Worksheets.Item("ModuleName")
$wsModule.Range("A3").Value = "Engagements"
$wsModule.Columns("A").ColumnWidth = 12.15

# --- AssociatedOpp sheet: move selection to A3 ---
$wsOpp = $wb.Worksheets.Item("AssociatedOpp")
$wsOpp.Range("A3").Select() | Out-Null

# --- AssociatedEng sheet: move selection to A2 ---
$wsEng = $wb.Worksheets.Item("AssociatedEng")
$wsEng.Range("A2").Select() | Out-Null

# --- Activate ModuleName last so it becomes the active/selected tab, with the selection on E21 ---
$wsModule.Activate() | Out-Null
$wsModule.Range("E21").Select() | Out-Null
